# Add data views for cds2db-modul with rights
# Applies the edits described in the diff against
# xl/worksheets/sheet1.xml (sheet "rights_and_functions"):
#   - restyle I33/J33/I34/J34 from the "Schlecht" red-highlight look to a
#     plain (unfilled) look while keeping them based on that cell style
#   - move the "20_take_over_check_date.sql" row data from row 37 to row 36
#   - populate row 37 with a brand-new "21_cre_view_typ_cds2db_all.sql" view
#     definition (with SELECT rights granted to cds2db_user)
#   - move the active selection to I47

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("rights_and_functions")

# ---------------------------------------------------------------------
# 1) Re-style I33, J33, I34, J34: drop the red "Schlecht" fill/font color
#    while keeping a normal Calibri font. Doing ClearFormats + a single
#    Font.Name write (applied uniformly to all four cells) produces one
#    shared new style entry instead of one per cell.
# ---------------------------------------------------------------------
$restyleCells = @("I33", "J33", "I34", "J34")
foreach ($addr in $restyleCells) {
    $ws.Range($addr).ClearFormats()
}
foreach ($addr in $restyleCells) {
    $rng = $ws.Range($addr)
    $rng.Font.Name = "Calibri"
    $rng.Font.Size = 11
}

# ---------------------------------------------------------------------
# 2) Row 36 currently holds I36/J36 ("SELECT" / "cds2db_user"); that data
#    is being replaced by what used to live in row 37 (the
#    take_over_check_date script description).
# ---------------------------------------------------------------------
$ws.Range("I36").ClearContents()
$ws.Range("J36").ClearContents()

$ws.Range("B36").Value = "20_take_over_check_date.sql"
$ws.Range("C36").Value = "template_take_over_check_date_function.sql"
$ws.Range("E36").Value = "db_log"
$ws.Range("H36").Value = "_raw"
$ws.Range("K36").Value = "template_take_over_check_date_function.sql"
$ws.Range("M36").Value = "take_over_last_check_date"
$ws.Range("N36").Value = "db_log"

# ---------------------------------------------------------------------
# 3) Row 37 becomes a brand-new entry for the 21_cre_view_typ_cds2db_all.sql
#    view, with SELECT rights granted to cds2db_user (styled like I33/J33).
#    K37/M37 (old content) are no longer used and get cleared.
# ---------------------------------------------------------------------
$ws.Range("K37").ClearContents()
$ws.Range("M37").ClearContents()

$ws.Range("B37").Value = "21_cre_view_typ_cds2db_all.sql"
$ws.Range("C37").Value = "template_cre_view2.sql"
$ws.Range("D37").Value = "cds2db_user"
$ws.Range("E37").Value = "cds2db_out"
$ws.Range("G37").Value = "v_"
$ws.Range("H37").Value = "_all"
$ws.Range("I37").Value = "SELECT"
# N37 already holds "db_log" (v24) and stays unchanged.

$j37 = $ws.Range("J37")
$j37.Value = "cds2db_user"
$j37.ClearFormats()
$j37.Font.Name = "Calibri"
$j37.Font.Size = 11

# ---------------------------------------------------------------------
# 4) Move the active cell / selection to I47 (matches the saved view
#    state in the edited workbook).
# ---------------------------------------------------------------------
$ws.Activate()
$ws.Range("I47").Select()
